$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Global")

# 1. Add the new "Aciclovir" worksheet right after "Global"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Aciclovir"

# 2. Move the Aciclovir parameter row (previously Global!A2:D2) onto the new sheet,
#    together with the shared header row, preserving number formatting.
$ws1.Range("A1:D2").Copy($ws2.Range("A1")) | Out-Null

# 3. Replace Global!A2:D2 with the new "Organism|Liver" / "EHC continuous fraction" row
$ws1.Range("B2").Value = "EHC continuous fraction"
$ws1.Range("A2").Value = "Organism|Liver"
$ws1.Range("C2").Value = 1
$ws1.Range("D2").ClearContents()

# 4. Update selections to match the saved view state
$ws2.Range("A1:XFD2").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C10").Select() | Out-Null
